$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (FIGUEROA GARCIA PIERINA GISELLA)
$ws.Range("E9").Value = 8
$ws.Range("H9").Value = 60

# Row 10 (LOZA GAVILANES ANDREA MAGDALENA)
$ws.Range("E10").Value = 8
$ws.Range("H10").Value = 60

# Row 16 (DELGADO LALANGUI LUIS GONZALO)
$ws.Range("E16").Value = 8
$ws.Range("H16").Value = 54

# Row 17 (ORRALA RIVAS JARITZA JAVIERA)
$ws.Range("E17").Value = 8
$ws.Range("H17").Value = 54

# Row 23 (ANDRADE VALVERDE PATRICIO SANTIAGO)
$ws.Range("E23").Value = 0
$ws.Range("H23").Value = 0

# Row 24 (CABRERA MOREIRA JURGEN WILLIAM)
$ws.Range("E24").Value = 3
$ws.Range("H24").Value = 49

# Update the selected cell shown in the bottom-right pane of the frozen view
$ws.Range("H17").Select() | Out-Null
